$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Write the new column Z data (values first, so dependent COUNTA/COUNTIF
#        formulas in columns B:I recalc correctly) -----------------------------

# Row 1: new training date (08/08/2025 -> serial 45877)
$ws.Range("Z1").Value = 45877

# Rows 2-27: attendance status for the new date.
# Default is "P" (Présent); a few players have a different status that day.
$attendance = @{
  2  = "P"
  3  = "P"
  4  = "P"
  5  = "P"
  6  = "P"
  7  = "P"
  8  = "P"
  9  = "P"
  10 = "P"
  11 = "P"
  12 = "P"
  13 = "P"
  14 = "P"
  15 = "P"
  16 = "P"
  17 = "P"
  18 = "P"
  19 = "P"
  20 = "P"
  21 = "P"
  22 = "P"
  23 = "P"
  24 = "P"
  25 = "M"
  26 = "P"
  27 = "RH"
}

foreach ($r in $attendance.Keys) {
  $ws.Range("Z$r").Value = $attendance[$r]
}

# --- 2. Copy the formatting from column Y onto the new column Z ---------------
$ws.Range("Y1:Y27").Copy()
$ws.Range("Z1:Z27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3. Restore the selection like the saved workbook (AB25) ------------------
$ws.Range("AB25").Select() | Out-Null
